$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E is "duplicate_image_filename". Fill "NA" for the data rows
# (the practice rows 2-5 and the main stimuli rows 6-21) that didn't
# already have a value in that column.
$ws.Range("E2:E21").Value = "NA"
